$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New task row data ("Realiseren programma", week 11-06 / 17-06) ---
# Row 7: finish the existing "Technisch ontwerp" block by adding its task + hours,
# then start a new block for "Realiseren programma".
$ws.Range("B7").Value = "Realiseren programma"
$ws.Range("C7").Value = 1.4583333333333333

$ws.Range("A8").Value = "11-06 / 17-06"
$ws.Range("B8").Value = "Realiseren programma"
$ws.Range("C8").Value = 0.16666666666666666

# --- Correction to an existing hour entry ---
$ws.Range("C3").Value = 0.020833333333333332

# --- Re-apply the elapsed-time ("Komma"/[h]:mm:ss) number format to the whole
#     hours column so it matches the single format already used by the total ---
$ws.Range("C2:C15").NumberFormat = "[h]:mm:ss"
$ws.Range("C16").NumberFormat = "[h]:mm:ss"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moved to F9 ---
$null = $ws.Range("F9").Select()
